$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.647.11'
$ws.Range("E2").Value = '  -2.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.893.96'
$ws.Range("E3").Value = '  -4.11%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.49'
$ws.Range("E5").Value = '  -2.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.14'
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("E8").Value = '  -3.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.895.11'
$ws.Range("E9").Value = '  -4.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  +4.79%  '

$ws.Range("E11").Value = '  -4.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.445'
$ws.Range("E12").Value = '  -3.11%  '

$ws.Range("E13").Value = '  -4.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.76'
$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.379.09'
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.603.58'
$ws.Range("E17").Value = '  -2.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.76'
$ws.Range("E18").Value = '  -3.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.894.83'
$ws.Range("E19").Value = '  -4.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '424.04'
$ws.Range("E20").Value = '  -5.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  -4.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.665'
$ws.Range("E22").Value = '  -3.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.04'
$ws.Range("E23").Value = '  -5.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.97'
$ws.Range("E24").Value = '  -2.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.80'
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.18'
$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.17'
$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("E32").Value = '  -3.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.27'
$ws.Range("E33").Value = '  -4.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  -3.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0827'
$ws.Range("E35").Value = '  -2.26%  '

$ws.Range("E36").Value = '  -2.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.61'
$ws.Range("E37").Value = '  -3.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.19'
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("E39").Value = '  -2.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  -2.84%  '

$ws.Range("E41").Value = '  +0.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.66'
$ws.Range("E42").Value = '  -4.06%  '

$ws.Range("E43").Value = '  +1.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.54'
$ws.Range("E44").Value = '  +3.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0344'
$ws.Range("E45").Value = '  -2.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '370.62'
$ws.Range("E46").Value = '  -5.27%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.647.81'
$ws.Range("E47").Value = '  -3.93%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.03'
$ws.Range("E48").Value = '  -1.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.67'
$ws.Range("E50").Value = '  +3.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.105'
$ws.Range("E51").Value = '  -1.89%  '
